$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12 (the "HEALSEC 20MG 14 CAPS" row), pushing that
# row and everything below it down by one. This keeps the alphabetically
# sorted product list in order and inserts the new product
# "GENICA RAY SUN SCREEN SPF 50+ CREAM 100 GM" between "EPIMAG EFFERVESCENT 12
# SACHETS" (row 11) and "HEALSEC 20MG 14 CAPS" (now row 13).
$ws.Rows.Item(12).Insert()

# Fill in the new row's data, mirroring the layout used by every other
# product row (A=index, B=name, H=ratio code, L=value, N=ratio code).
$ws.Cells.Item(12, 1).Value = 9
$ws.Cells.Item(12, 2).Value = "GENICA RAY SUN SCREEN SPF 50+ CREAM 100 GM"
$ws.Cells.Item(12, 8).Value = "0:0"
$ws.Cells.Item(12, 12).Value = 198
$ws.Cells.Item(12, 14).Value = "1:0"

# Renumber column A for every row below the inserted one (10, 11, 12, ...).
for ($r = 13; $r -le 27; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 3
}
